$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-13 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-14 Saturday", 2) | Out-Null
$d.Content.Find.Execute("10+58=", $true, $false, $false, $false, $false, $true, 1, $false, "93-65=", 2) | Out-Null
$d.Content.Find.Execute("53+24=", $true, $false, $false, $false, $false, $true, 1, $false, "87+7=", 2) | Out-Null
$d.Content.Find.Execute("14+3=", $true, $false, $false, $false, $false, $true, 1, $false, "55-46=", 2) | Out-Null
$d.Content.Find.Execute("74-68=", $true, $false, $false, $false, $false, $true, 1, $false, "51+10=", 2) | Out-Null
$d.Content.Find.Execute("20+9=", $true, $false, $false, $false, $false, $true, 1, $false, "58-11=", 2) | Out-Null
$d.Content.Find.Execute("50-26=", $true, $false, $false, $false, $false, $true, 1, $false, "77-66=", 2) | Out-Null
$d.Content.Find.Execute("25+0=", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=", 2) | Out-Null
$d.Content.Find.Execute("81-79=", $true, $false, $false, $false, $false, $true, 1, $false, "56-29=", 2) | Out-Null
$d.Content.Find.Execute("79+7=", $true, $false, $false, $false, $false, $true, 1, $false, "87-6=", 2) | Out-Null
$d.Content.Find.Execute("37-30=", $true, $false, $false, $false, $false, $true, 1, $false, "18-15=", 2) | Out-Null
$d.Content.Find.Execute("24-1=", $true, $false, $false, $false, $false, $true, 1, $false, "28+46=", 2) | Out-Null
$d.Content.Find.Execute("57-34=", $true, $false, $false, $false, $false, $true, 1, $false, "78-6=", 2) | Out-Null
$d.Content.Find.Execute("67+20=", $true, $false, $false, $false, $false, $true, 1, $false, "72-54=", 2) | Out-Null
$d.Content.Find.Execute("50+5=", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=", 2) | Out-Null
$d.Content.Find.Execute("0+6=", $true, $false, $false, $false, $false, $true, 1, $false, "76-11=", 2) | Out-Null
$d.Content.Find.Execute("71+8=", $true, $false, $false, $false, $false, $true, 1, $false, "58-4=", 2) | Out-Null
$d.Content.Find.Execute("68+8=", $true, $false, $false, $false, $false, $true, 1, $false, "16+79=", 2) | Out-Null
$d.Content.Find.Execute("13+74=", $true, $false, $false, $false, $false, $true, 1, $false, "63-29=", 2) | Out-Null
$d.Content.Find.Execute("16-13=", $true, $false, $false, $false, $false, $true, 1, $false, "66-1=", 2) | Out-Null
$d.Content.Find.Execute("85+5=", $true, $false, $false, $false, $false, $true, 1, $false, "89-88=", 2) | Out-Null
$d.Content.Find.Execute("64+1=", $true, $false, $false, $false, $false, $true, 1, $false, "97-26=", 2) | Out-Null
$d.Content.Find.Execute("68+26=", $true, $false, $false, $false, $false, $true, 1, $false, "23+56=", 2) | Out-Null
$d.Content.Find.Execute("83+9=", $true, $false, $false, $false, $false, $true, 1, $false, "5+46=", 2) | Out-Null
$d.Content.Find.Execute("62-44=", $true, $false, $false, $false, $false, $true, 1, $false, "68-36=", 2) | Out-Null
$d.Content.Find.Execute("26+47=", $true, $false, $false, $false, $false, $true, 1, $false, "12+45=", 2) | Out-Null
$d.Content.Find.Execute("35+40=", $true, $false, $false, $false, $false, $true, 1, $false, "82-63=", 2) | Out-Null
$d.Content.Find.Execute("47-10=", $true, $false, $false, $false, $false, $true, 1, $false, "75-55=", 2) | Out-Null
$d.Content.Find.Execute("28+15=", $true, $false, $false, $false, $false, $true, 1, $false, "51+35=", 2) | Out-Null
$d.Content.Find.Execute("16+4=", $true, $false, $false, $false, $false, $true, 1, $false, "43+24=", 2) | Out-Null
$d.Content.Find.Execute("17+59=", $true, $false, $false, $false, $false, $true, 1, $false, "41+12=", 2) | Out-Null
$d.Content.Find.Execute("53+31=", $true, $false, $false, $false, $false, $true, 1, $false, "20+27=", 2) | Out-Null
$d.Content.Find.Execute("47-42=", $true, $false, $false, $false, $false, $true, 1, $false, "44+37=", 2) | Out-Null
$d.Content.Find.Execute("87-11=", $true, $false, $false, $false, $false, $true, 1, $false, "97-75=", 2) | Out-Null
$d.Content.Find.Execute("57-47=", $true, $false, $false, $false, $false, $true, 1, $false, "0+85=", 2) | Out-Null
$d.Content.Find.Execute("64-17=", $true, $false, $false, $false, $false, $true, 1, $false, "38+1=", 2) | Out-Null
$d.Content.Find.Execute("36-8=", $true, $false, $false, $false, $false, $true, 1, $false, "78-46=", 2) | Out-Null
$d.Content.Find.Execute("43-17=", $true, $false, $false, $false, $false, $true, 1, $false, "85-85=", 2) | Out-Null
$d.Content.Find.Execute("95-10=", $true, $false, $false, $false, $false, $true, 1, $false, "39+33=", 2) | Out-Null
$d.Content.Find.Execute("98-42=", $true, $false, $false, $false, $false, $true, 1, $false, "49+38=", 2) | Out-Null
$d.Content.Find.Execute("81+8=", $true, $false, $false, $false, $false, $true, 1, $false, "24+4=", 2) | Out-Null
$d.Content.Find.Execute("32-3=", $true, $false, $false, $false, $false, $true, 1, $false, "8+10=", 2) | Out-Null
$d.Content.Find.Execute("69-54=", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=", 2) | Out-Null
$d.Content.Find.Execute("86-65=", $true, $false, $false, $false, $false, $true, 1, $false, "82-78=", 2) | Out-Null
$d.Content.Find.Execute("50-28=", $true, $false, $false, $false, $false, $true, 1, $false, "22+11=", 2) | Out-Null
$d.Content.Find.Execute("0+22=", $true, $false, $false, $false, $false, $true, 1, $false, "22+37=", 2) | Out-Null
$d.Content.Find.Execute("66+17=", $true, $false, $false, $false, $false, $true, 1, $false, "2+68=", 2) | Out-Null
$d.Content.Find.Execute("69-23=", $true, $false, $false, $false, $false, $true, 1, $false, "78-76=", 2) | Out-Null
$d.Content.Find.Execute("17-15=", $true, $false, $false, $false, $false, $true, 1, $false, "75-69=", 2) | Out-Null
$d.Content.Find.Execute("75-41=", $true, $false, $false, $false, $false, $true, 1, $false, "42-37=", 2) | Out-Null
$d.Content.Find.Execute("34+60=", $true, $false, $false, $false, $false, $true, 1, $false, "31+9=", 2) | Out-Null
$d.Content.Find.Execute("19+38=", $true, $false, $false, $false, $false, $true, 1, $false, "75-4=", 2) | Out-Null
$d.Content.Find.Execute("56-22=", $true, $false, $false, $false, $false, $true, 1, $false, "75+10=", 2) | Out-Null
$d.Content.Find.Execute("0+5=", $true, $false, $false, $false, $false, $true, 1, $false, "80-62=", 2) | Out-Null
$d.Content.Find.Execute("5+16=", $true, $false, $false, $false, $false, $true, 1, $false, "62-28=", 2) | Out-Null
$d.Content.Find.Execute("76-36=", $true, $false, $false, $false, $false, $true, 1, $false, "67-21=", 2) | Out-Null
$d.Content.Find.Execute("6+34=", $true, $false, $false, $false, $false, $true, 1, $false, "18-12=", 2) | Out-Null
$d.Content.Find.Execute("25-0=", $true, $false, $false, $false, $false, $true, 1, $false, "45+25=", 2) | Out-Null
$d.Content.Find.Execute("5+11=", $true, $false, $false, $false, $false, $true, 1, $false, "61+20=", 2) | Out-Null
$d.Content.Find.Execute("16+71=", $true, $false, $false, $false, $false, $true, 1, $false, "64-56=", 2) | Out-Null
$d.Content.Find.Execute("76+2=", $true, $false, $false, $false, $false, $true, 1, $false, "9+12=", 2) | Out-Null
$d.Content.Find.Execute("68-63=", $true, $false, $false, $false, $false, $true, 1, $false, "31+28=", 2) | Out-Null
$d.Content.Find.Execute("45+37=", $true, $false, $false, $false, $false, $true, 1, $false, "30-16=", 2) | Out-Null
$d.Content.Find.Execute("86-84=", $true, $false, $false, $false, $false, $true, 1, $false, "84+10=", 2) | Out-Null
$d.Content.Find.Execute("31-31=", $true, $false, $false, $false, $false, $true, 1, $false, "93-53=", 2) | Out-Null
$d.Content.Find.Execute("81-5=", $true, $false, $false, $false, $false, $true, 1, $false, "8+18=", 2) | Out-Null
$d.Content.Find.Execute("53+10=", $true, $false, $false, $false, $false, $true, 1, $false, "83-46=", 2) | Out-Null
$d.Content.Find.Execute("54-6=", $true, $false, $false, $false, $false, $true, 1, $false, "52+7=", 2) | Out-Null
$d.Content.Find.Execute("3+21=", $true, $false, $false, $false, $false, $true, 1, $false, "95-65=", 2) | Out-Null
$d.Content.Find.Execute("75-8=", $true, $false, $false, $false, $false, $true, 1, $false, "80-77=", 2) | Out-Null
$d.Content.Find.Execute("71-5=", $true, $false, $false, $false, $false, $true, 1, $false, "79-66=", 2) | Out-Null
$d.Content.Find.Execute("18+22=", $true, $false, $false, $false, $false, $true, 1, $false, "45-26=", 2) | Out-Null
$d.Content.Find.Execute("65-18=", $true, $false, $false, $false, $false, $true, 1, $false, "6+14=", 2) | Out-Null
$d.Content.Find.Execute("56+41=", $true, $false, $false, $false, $false, $true, 1, $false, "61+11=", 2) | Out-Null
$d.Content.Find.Execute("90-61=", $true, $false, $false, $false, $false, $true, 1, $false, "84-59=", 2) | Out-Null
$d.Content.Find.Execute("35-25=", $true, $false, $false, $false, $false, $true, 1, $false, "10+79=", 2) | Out-Null
$d.Content.Find.Execute("40-14=", $true, $false, $false, $false, $false, $true, 1, $false, "73-18=", 2) | Out-Null
$d.Content.Find.Execute("59+39=", $true, $false, $false, $false, $false, $true, 1, $false, "28+59=", 2) | Out-Null
$d.Content.Find.Execute("80-49=", $true, $false, $false, $false, $false, $true, 1, $false, "2+97=", 2) | Out-Null
$d.Content.Find.Execute("84-33=", $true, $false, $false, $false, $false, $true, 1, $false, "91-31=", 2) | Out-Null
$d.Content.Find.Execute("24+31=", $true, $false, $false, $false, $false, $true, 1, $false, "92-89=", 2) | Out-Null
$d.Content.Find.Execute("57-43=", $true, $false, $false, $false, $false, $true, 1, $false, "39+51=", 2) | Out-Null
$d.Content.Find.Execute("16+49=", $true, $false, $false, $false, $false, $true, 1, $false, "32-4=", 2) | Out-Null
$d.Content.Find.Execute("32+45=", $true, $false, $false, $false, $false, $true, 1, $false, "74-66=", 2) | Out-Null
$d.Content.Find.Execute("99-34=", $true, $false, $false, $false, $false, $true, 1, $false, "95-46=", 2) | Out-Null
$d.Content.Find.Execute("27+60=", $true, $false, $false, $false, $false, $true, 1, $false, "78-56=", 2) | Out-Null
$d.Content.Find.Execute("82+16=", $true, $false, $false, $false, $false, $true, 1, $false, "32-10=", 2) | Out-Null
$d.Content.Find.Execute("28-1=", $true, $false, $false, $false, $false, $true, 1, $false, "5+39=", 2) | Out-Null
$d.Content.Find.Execute("35-35=", $true, $false, $false, $false, $false, $true, 1, $false, "69-3=", 2) | Out-Null
$d.Content.Find.Execute("72-22=", $true, $false, $false, $false, $false, $true, 1, $false, "32+46=", 2) | Out-Null
$d.Content.Find.Execute("46+31=", $true, $false, $false, $false, $false, $true, 1, $false, "76-58=", 2) | Out-Null
$d.Content.Find.Execute("83-56=", $true, $false, $false, $false, $false, $true, 1, $false, "36+38=", 2) | Out-Null
$d.Content.Find.Execute("10+81=", $true, $false, $false, $false, $false, $true, 1, $false, "67-42=", 2) | Out-Null
$d.Content.Find.Execute("5+87=", $true, $false, $false, $false, $false, $true, 1, $false, "28-16=", 2) | Out-Null
$d.Content.Find.Execute("53-42=", $true, $false, $false, $false, $false, $true, 1, $false, "86-52=", 2) | Out-Null
$d.Content.Find.Execute("77-60=", $true, $false, $false, $false, $false, $true, 1, $false, "74-10=", 2) | Out-Null
$d.Content.Find.Execute("57-36=", $true, $false, $false, $false, $false, $true, 1, $false, "66-31=", 2) | Out-Null
$d.Content.Find.Execute("65-56=", $true, $false, $false, $false, $false, $true, 1, $false, "52-5=", 2) | Out-Null
$d.Content.Find.Execute("11+86=", $true, $false, $false, $false, $false, $true, 1, $false, "74-15=", 2) | Out-Null
$d.Content.Find.Execute("63-14=", $true, $false, $false, $false, $false, $true, 1, $false, "48-0=", 2) | Out-Null
$d.Content.Find.Execute("10+85=", $true, $false, $false, $false, $false, $true, 1, $false, "34+64=", 2) | Out-Null
